$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the "last updated" timestamp string shown in A1
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 16 de Septiembre de 2020 a las 05:36"

# ---------------------------------------------------------------------------
# 2. Re-rank countries whose position in the table changed.
#    Country names live in column A; use a unique temporary placeholder while
#    swapping two names so that the two operations don't collide when one of
#    the target names already exists elsewhere in the sheet.
# ---------------------------------------------------------------------------

# Barein (row 53) <-> Venezuela (row 54)
$ws.Range("A53").Value = "__TMP_SWAP_A__"
$ws.Range("A54").Value = "Barein"
$ws.Range("A53").Value = "Venezuela"

# Siria (row 128) -> Birmania, Angola (row 129) -> Siria, Birmania (row 130) -> Angola
$ws.Range("A128").Value = "__TMP_SWAP_B__"
$ws.Range("A129").Value = "__TMP_SWAP_C__"
$ws.Range("A130").Value = "Angola"
$ws.Range("A129").Value = "Siria"
$ws.Range("A128").Value = "Birmania"

# Timor Oriental (row 204) <-> Santa Lucia (row 205)
$ws.Range("A204").Value = "__TMP_SWAP_D__"
$ws.Range("A205").Value = "Timor Oriental"
$ws.Range("A204").Value = "Santa Lucia"

# Islas Malvinas (row 214) <-> Montserrat (row 215)
$ws.Range("A214").Value = "__TMP_SWAP_E__"
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("A214").Value = "Montserrat"

# ---------------------------------------------------------------------------
# 3. Refresh the statistics (Casos totales, Nuevos casos, Casos activos,
#    Recuperados, Casos criticos, Muertes hoy, Muertes) for every row whose
#    figures moved.
# ---------------------------------------------------------------------------

# Row 20 - Pakistan
$ws.Range("B20").Value = 303089
$ws.Range("C20").Value = 665
$ws.Range("D20").Value = 290760
$ws.Range("E20").Value = 5936
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 6393

# Row 33 - Kazajistan
$ws.Range("B33").Value = 106984
$ws.Range("C33").Value = 64
$ws.Range("D33").Value = 101267
$ws.Range("E33").Value = 4046
$ws.Range("H33").Value = 1671

# Row 39 - Belgica
$ws.Range("B39").Value = 94795
$ws.Range("C39").Value = 489
$ws.Range("D39").Value = 18789
$ws.Range("E39").Value = 66076
$ws.Range("G39").Value = 3
$ws.Range("H39").Value = 9930

# Row 53 - now Venezuela
$ws.Range("B53").Value = 62655
$ws.Range("D53").Value = 50361
$ws.Range("E53").Value = 11792
$ws.Range("H53").Value = 502

# Row 54 - now Barein
$ws.Range("B54").Value = 61643
$ws.Range("D54").Value = 54831
$ws.Range("E54").Value = 6599
$ws.Range("H54").Value = 213

# Row 128 - now Birmania
$ws.Range("B128").Value = 3636
$ws.Range("C128").Value = 134
$ws.Range("D128").Value = 832
$ws.Range("E128").Value = 2765
$ws.Range("G128").Value = 4
$ws.Range("H128").Value = 39

# Row 129 - now Siria
$ws.Range("B129").Value = 3614
$ws.Range("D129").Value = 871
$ws.Range("E129").Value = 2583
$ws.Range("H129").Value = 160

# Row 130 - now Angola
$ws.Range("B130").Value = 3569
$ws.Range("D130").Value = 1332
$ws.Range("E130").Value = 2098
$ws.Range("H130").Value = 139

# Row 172 - Islas Turcas y Caicos
$ws.Range("B172").Value = 650
$ws.Range("C172").Value = 2
$ws.Range("D172").Value = 557
$ws.Range("E172").Value = 88

# Row 173 - San Martin (Parte Holandesa)
$ws.Range("D173").Value = 468
$ws.Range("E173").Value = 62

# Row 214 - now Montserrat
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

# Row 215 - now Islas Malvinas
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
